$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17 ---
$ws.Range("A17").Value = 45204
$ws.Range("A2").Copy()
$ws.Range("A17").PasteSpecial(-4122)
$ws.Range("B17").Value = "MPAL"
$ws.Range("C17").Value = "TP"
$ws.Range("D17").Value = "X"
$ws.Rows.Item(17).RowHeight = 232.5
$ws.Range("H17").Value = 'Reformulation de la question elaastic 3 ?'
$ws.Range("G17").Value = 'Elaastic : question 2 et 3 : difficulté du groupe à comprendre la signification de la question 3 "comment aurait pu être traitée" ne leur semble pas très clair.
Suite des rédactions des Tests d''Acceptation : remarques fréquentes sur des éléments du When qui pourraient en fait être placés dans le Given : "Given : Un membre non connecté sur la page de connexion
When : il clique sur le bouton de connexion après avoir correctement rempli ses informations
Then : il est connecté"
Le "après avoir correctement rempli ses informations" devrait être placé dans le Given ?'

# --- Row 18 ---
$ws.Range("A18").Value = 45204
$ws.Range("A2").Copy()
$ws.Range("A18").PasteSpecial(-4122)
$ws.Range("B18").Value = "MPAL"
$ws.Range("C18").Value = "TP"
$ws.Rows.Item(18).RowHeight = 217
$ws.Range("F18").Value = "X"
$ws.Range("G18").Value = 'Elaastic : fin de l''activité + Rédaction de Tests d''acceptation. Pourrait-on envisager d''avoir des Business Rules et Lexique externe à toutes US pour préciser le type de visiteur ? Exemple :
"Business Rule : types d''utilisateur
Il y a 4 types d''utilisateurs du site ALOSA :
- visiteur : un internaute non authentifié
- membre : internaute authentifié bénéficiant de permissions standard
- membre expert : internaute authentifié bénéficiant de permissions expert
- administrateur : internaute authentifié bénéficiant de permissions administrateur
"'

$excel.CutCopyMode = 0

# --- View state ---
$null = $ws.Range("G19").Select()
